$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 - 2023-08-10 morning reading
$ws.Range("A8").Value = 45148.472916666666
$ws.Range("B8").Formula = "=(123+126)/2"
$ws.Range("C8").Formula = "=(81+85)/2"
$ws.Range("D8").Formula = "=(95+89)/2"
$ws.Range("E8").Value = 97

# Row 9 - 2023-08-10 afternoon reading
$ws.Range("A9").Value = 45148.597222222219
$ws.Range("B9").Formula = "=(130+124)/2"
$ws.Range("C9").Formula = "=(90+88)/2"
$ws.Range("D9").Formula = "=(89+86)/2"
$ws.Range("E9").Value = 98

# Row 10 - 2023-08-10 evening reading
$ws.Range("A10").Value = 45148.780555555553
$ws.Range("B10").Formula = "=(120+129)/2"
$ws.Range("C10").Formula = "=(90+89)/2"
$ws.Range("D10").Formula = "=(115+107)/2"
$ws.Range("E10").Value = 95

# Update the active selection to A11, matching the saved cursor position
[void]$ws.Range("A11").Select()
